# Helper: write a cell as genuine text (avoids Excel's automatic
# "looks like a number" coercion that would otherwise turn "011531" into
# 11531, or flip t="inlineStr" into t="n"). We stamp the cell's number
# format to Text ("@") before assigning the string, then restore the
# cell's formatting back to a neighboring plain/default-styled cell so we
# do not leave a stray derived style behind.
function Set-TextValue($sheet, $row, $col, $text) {
    $cell = $sheet.Cells.Item($row, $col)
    $blank = $sheet.Cells.Item(200, 200)
    $blank.Copy()
    $cell.PasteSpecial(-4122)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $blank.Copy()
    $cell.PasteSpecial(-4122)
}

$wb = $excel.ActiveWorkbook

# --- Step 1: structural changes first -----------------------------------
# Insert the new "2022-Q1" sheet right after "2021-Q4" (i.e. right before
# "总计"), mirroring the 2021-Q4 fund-holdings layout. Sheet
# references are re-fetched (by name) AFTER this call below, since adding
# a sheet can shift sheet positions and invalidate previously-held
# references.
$q4ForAnchor = $wb.Worksheets.Item("2021-Q4")
$newSheet = $wb.Worksheets.Add($null, $q4ForAnchor)
$newSheet.Name = "2022-Q1"

# --- Step 2: re-fetch every sheet reference now that the sheet list is
# final, then do all the cell writes. ------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$total = $wb.Worksheets.Item("总计")

# Clone the 2021-Q4 header/column formatting (bold+border on the header
# row and the index column) onto the new sheet before filling in values.
# (Copied as two ranges, deliberately skipping A1 - it is blank/unused in
# 2021-Q4 too, and we do not want to materialize a stray empty cell there.)
$q4.Range("B1:H1").Copy()
$newSheet.Range("B1").PasteSpecial(-4122)
$q4.Range("A2:A14").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)

$newSheet.Cells.Item(1,2).Value = "基金代码"
$newSheet.Cells.Item(1,3).Value = "基金名称"
$newSheet.Cells.Item(1,4).Value = "基金规模"
$newSheet.Cells.Item(1,5).Value = "股票总仓位"
$newSheet.Cells.Item(1,6).Value = "仓位占比"
$newSheet.Cells.Item(1,7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1,8).Value = "仓位排名"
$newSheet.Cells.Item(2,1).Value = 0
Set-TextValue $newSheet 2 2 "011531"
Set-TextValue $newSheet 2 3 "朱雀恒心一年持有期混合"
Set-TextValue $newSheet 2 4 "90.24"
Set-TextValue $newSheet 2 5 "93.49"
Set-TextValue $newSheet 2 6 "6.55"
Set-TextValue $newSheet 2 7 "5.9107"
$newSheet.Cells.Item(2,8).Value = 3
$newSheet.Cells.Item(3,1).Value = 1
Set-TextValue $newSheet 3 2 "007493"
Set-TextValue $newSheet 3 3 "朱雀产业臻选混合A"
Set-TextValue $newSheet 3 4 "51.18"
Set-TextValue $newSheet 3 5 "93.91"
Set-TextValue $newSheet 3 6 "5.24"
Set-TextValue $newSheet 3 7 "2.6818"
$newSheet.Cells.Item(3,8).Value = 4
$newSheet.Cells.Item(4,1).Value = 2
Set-TextValue $newSheet 4 2 "010141"
Set-TextValue $newSheet 4 3 "朱雀企业优选股票A"
Set-TextValue $newSheet 4 4 "36.18"
Set-TextValue $newSheet 4 5 "92.72"
Set-TextValue $newSheet 4 6 "6.63"
Set-TextValue $newSheet 4 7 "2.3987"
$newSheet.Cells.Item(4,8).Value = 3
$newSheet.Cells.Item(5,1).Value = 3
Set-TextValue $newSheet 5 2 "010922"
Set-TextValue $newSheet 5 3 "朱雀匠心一年持有期混合"
Set-TextValue $newSheet 5 4 "25.58"
Set-TextValue $newSheet 5 5 "94.26"
Set-TextValue $newSheet 5 6 "9.02"
Set-TextValue $newSheet 5 7 "2.3073"
$newSheet.Cells.Item(5,8).Value = 1
$newSheet.Cells.Item(6,1).Value = 4
Set-TextValue $newSheet 6 2 "008294"
Set-TextValue $newSheet 6 3 "朱雀企业优胜股票A"
Set-TextValue $newSheet 6 4 "15.13"
Set-TextValue $newSheet 6 5 "93.46"
Set-TextValue $newSheet 6 6 "5.23"
Set-TextValue $newSheet 6 7 "0.7913"
$newSheet.Cells.Item(6,8).Value = 4
$newSheet.Cells.Item(7,1).Value = 5
Set-TextValue $newSheet 7 2 "007494"
Set-TextValue $newSheet 7 3 "朱雀产业臻选混合C"
Set-TextValue $newSheet 7 4 "13.56"
Set-TextValue $newSheet 7 5 "93.91"
Set-TextValue $newSheet 7 6 "5.24"
Set-TextValue $newSheet 7 7 "0.7105"
$newSheet.Cells.Item(7,8).Value = 4
$newSheet.Cells.Item(8,1).Value = 6
Set-TextValue $newSheet 8 2 "910005"
Set-TextValue $newSheet 8 3 "东方红启兴三年持有混合A"
Set-TextValue $newSheet 8 4 "6.24"
Set-TextValue $newSheet 8 5 "92.98"
Set-TextValue $newSheet 8 6 "6.56"
Set-TextValue $newSheet 8 7 "0.4093"
$newSheet.Cells.Item(8,8).Value = 6
$newSheet.Cells.Item(9,1).Value = 7
Set-TextValue $newSheet 9 2 "010142"
Set-TextValue $newSheet 9 3 "朱雀企业优选股票C"
Set-TextValue $newSheet 9 4 "5.30"
Set-TextValue $newSheet 9 5 "92.72"
Set-TextValue $newSheet 9 6 "6.63"
Set-TextValue $newSheet 9 7 "0.3514"
$newSheet.Cells.Item(9,8).Value = 3
$newSheet.Cells.Item(10,1).Value = 8
Set-TextValue $newSheet 10 2 "007880"
Set-TextValue $newSheet 10 3 "朱雀产业智选混合A"
Set-TextValue $newSheet 10 4 "5.06"
Set-TextValue $newSheet 10 5 "92.87"
Set-TextValue $newSheet 10 6 "5.38"
Set-TextValue $newSheet 10 7 "0.2722"
$newSheet.Cells.Item(10,8).Value = 3
$newSheet.Cells.Item(11,1).Value = 9
Set-TextValue $newSheet 11 2 "008295"
Set-TextValue $newSheet 11 3 "朱雀企业优胜股票C"
Set-TextValue $newSheet 11 4 "2.76"
Set-TextValue $newSheet 11 5 "93.46"
Set-TextValue $newSheet 11 6 "5.23"
Set-TextValue $newSheet 11 7 "0.1443"
$newSheet.Cells.Item(11,8).Value = 4
$newSheet.Cells.Item(12,1).Value = 10
Set-TextValue $newSheet 12 2 "007881"
Set-TextValue $newSheet 12 3 "朱雀产业智选混合C"
Set-TextValue $newSheet 12 4 "0.80"
Set-TextValue $newSheet 12 5 "92.87"
Set-TextValue $newSheet 12 6 "5.38"
Set-TextValue $newSheet 12 7 "0.0430"
$newSheet.Cells.Item(12,8).Value = 3
$newSheet.Cells.Item(13,1).Value = 11
Set-TextValue $newSheet 13 2 "580007"
Set-TextValue $newSheet 13 3 "东吴安享量化灵活配置混合"
Set-TextValue $newSheet 13 4 "0.57"
Set-TextValue $newSheet 13 5 "73.95"
Set-TextValue $newSheet 13 6 "4.97"
Set-TextValue $newSheet 13 7 "0.0283"
$newSheet.Cells.Item(13,8).Value = 10
$newSheet.Cells.Item(14,1).Value = 12
Set-TextValue $newSheet 14 2 "690003"
Set-TextValue $newSheet 14 3 "民生加银精选混合"
Set-TextValue $newSheet 14 4 "0.61"
Set-TextValue $newSheet 14 5 "91.83"
Set-TextValue $newSheet 14 6 "4.06"
Set-TextValue $newSheet 14 7 "0.0248"
$newSheet.Cells.Item(14,8).Value = 9

# 2) Prepend the 2022-Q1 summary row to the "总计" sheet, shifting
#    the existing rows down by one and renumbering the index column.

# Give the new last row (row 6) the same index-column style as the row
# above it before overwriting values top-to-bottom.
$total.Cells.Item(5, 1).Copy()
$total.Cells.Item(6, 1).PasteSpecial(-4122)

$total.Cells.Item(6, 1).Value = 4
$total.Cells.Item(6, 2).Value = "2021-Q1"
$total.Cells.Item(6, 3).Value = 2
$total.Cells.Item(6, 4).Value = 0.23

$total.Cells.Item(5, 1).Value = 3
$total.Cells.Item(5, 2).Value = "2021-Q2"
$total.Cells.Item(5, 3).Value = 6
$total.Cells.Item(5, 4).Value = 9.69

$total.Cells.Item(4, 1).Value = 2
$total.Cells.Item(4, 2).Value = "2021-Q3"
$total.Cells.Item(4, 3).Value = 10
$total.Cells.Item(4, 4).Value = 15.74

$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(3, 2).Value = "2021-Q4"
$total.Cells.Item(3, 3).Value = 13
$total.Cells.Item(3, 4).Value = 12.33

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q1"
$total.Cells.Item(2, 3).Value = 13
$total.Cells.Item(2, 4).Value = 16.07
